# Finished tweaks on data-viz-01
#
# The "Exercise, change the color" deck (slides 3-7) gets speaker notes
# added -- explaining the Altair/Python, R and Tableau outputs. Slide 2
# already carries a notes page (notesSlide1.xml); this adds the
# matching notes pages for slides 3 through 7.

$p = $ppt.ActivePresentation

function Set-SpeakerNotes {
    param(
        [int]$SlideIndex,
        [string[]]$Paragraphs
    )

    $slide = $p.Slides.Item($SlideIndex)
    $notesPage = $slide.NotesPage

    # This host only materializes the notes body placeholder on demand;
    # grab it (creating the real notesSlide part) before touching text.
    $notesShape = $notesPage.Shapes.AddPlaceholder(2)

    $nl = [char]10
    $notesShape.TextFrame.TextRange.Text = [string]::Join($nl, $Paragraphs)
}

# Slide 3 - "Exercise, Python code"
Set-SpeakerNotes 3 @(
    "In Altair/Python, the variable associated with different colors goes inside the encode function."
)

# Slide 4 - "Exercise, Python output"
Set-SpeakerNotes 4 @(
    "Notice that Altair/Python used a gradient of colors. This is a good choice.",
    "",
    "The darker colors, associated with a larger number of bathrooms cluster to the left. Newer houses tend to have more bathrooms."
)

# Slide 5 - "Exercise, R code"
Set-SpeakerNotes 5 @(
    "The color variable is defined inside the aes function."
)

# Slide 6 - "Exercise, R output"
Set-SpeakerNotes 6 @(
    "I deliberately defined Baths as categorical with the factor function, and notice the use of discrete colors again rather than a gradient."
)

# Slide 7 - "Exercise, Tableau output"
Set-SpeakerNotes 7 @(
    "I deliberately defined Baths as Dimension Categorical (blue pill) to show the contrast. Notice that the colors are discrete well separated values."
)
